$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I20").Value = 'ba'
$ws.Range("J20").Value = 'Appreciation'
$ws.Range("I21").Value = 'sv'
$ws.Range("J21").Value = 'Statement-opinion'
$ws.Range("I39").Value = 'sv'
$ws.Range("J39").Value = 'Statement-opinion'
$ws.Range("I42").Value = 'sd'
$ws.Range("J42").Value = 'Statement-non-opinion'
$ws.Range("I44").Value = 'sd'
$ws.Range("J44").Value = 'Statement-non-opinion'
$ws.Range("I49").Value = 'sv'
$ws.Range("J49").Value = 'Statement-opinion'
$ws.Range("I51").Value = 'sv'
$ws.Range("J51").Value = 'Statement-opinion'
$ws.Range("I53").Value = 'aa'
$ws.Range("J53").Value = 'Agree/Accept'
$ws.Range("I62").Value = 'ba'
$ws.Range("J62").Value = 'Appreciation'
$ws.Range("I67").Value = 'b'
$ws.Range("J67").Value = 'Acknowledge (Backchannel)'
$ws.Range("I88").Value = 'sd'
$ws.Range("J88").Value = 'Statement-non-opinion'
$ws.Range("I91").Value = 'b'
$ws.Range("J91").Value = 'Acknowledge (Backchannel)'
$ws.Range("I99").Value = 'sv'
$ws.Range("J99").Value = 'Statement-opinion'
$ws.Range("I111").Value = 'sv'
$ws.Range("J111").Value = 'Statement-opinion'
$ws.Range("I116").Value = 'sv'
$ws.Range("J116").Value = 'Statement-opinion'
$ws.Range("I117").Value = 'aa'
$ws.Range("J117").Value = 'Agree/Accept'
$ws.Range("I118").Value = 'b'
$ws.Range("J118").Value = 'Acknowledge (Backchannel)'
$ws.Range("I128").Value = 'sd'
$ws.Range("J128").Value = 'Statement-non-opinion'
$ws.Range("I134").Value = 'b'
$ws.Range("J134").Value = 'Acknowledge (Backchannel)'
$ws.Range("I135").Value = 'b'
$ws.Range("J135").Value = 'Acknowledge (Backchannel)'
$ws.Range("I143").Value = 'sd'
$ws.Range("J143").Value = 'Statement-non-opinion'
$ws.Range("I148").Value = 'sd'
$ws.Range("J148").Value = 'Statement-non-opinion'
$ws.Range("I167").Value = 'sd'
$ws.Range("J167").Value = 'Statement-non-opinion'
$ws.Range("I179").Value = 'b'
$ws.Range("J179").Value = 'Acknowledge (Backchannel)'
$ws.Range("I183").Value = 'b'
$ws.Range("J183").Value = 'Acknowledge (Backchannel)'
$ws.Range("I207").Value = 'ba'
$ws.Range("J207").Value = 'Appreciation'
$ws.Range("I243").Value = 'b'
$ws.Range("J243").Value = 'Acknowledge (Backchannel)'
$ws.Range("I255").Value = 'aa'
$ws.Range("J255").Value = 'Agree/Accept'
$ws.Range("I278").Value = 'b'
$ws.Range("J278").Value = 'Acknowledge (Backchannel)'
$ws.Range("I293").Value = 'sv'
$ws.Range("J293").Value = 'Statement-opinion'
$ws.Range("I322").Value = 'sv'
$ws.Range("J322").Value = 'Statement-opinion'
$ws.Range("I325").Value = 'sd'
$ws.Range("J325").Value = 'Statement-non-opinion'
$ws.Range("I336").Value = 'sv'
$ws.Range("J336").Value = 'Statement-opinion'
$ws.Range("I342").Value = 'b'
$ws.Range("J342").Value = 'Acknowledge (Backchannel)'
$ws.Range("I349").Value = 'b'
$ws.Range("J349").Value = 'Acknowledge (Backchannel)'
$ws.Range("I379").Value = 'b'
$ws.Range("J379").Value = 'Acknowledge (Backchannel)'
$ws.Range("I381").Value = 'sv'
$ws.Range("J381").Value = 'Statement-opinion'
$ws.Range("I384").Value = 'sv'
$ws.Range("J384").Value = 'Statement-opinion'
$ws.Range("I394").Value = 'ba'
$ws.Range("J394").Value = 'Appreciation'
$ws.Range("I407").Value = 'sv'
$ws.Range("J407").Value = 'Statement-opinion'
$ws.Range("I410").Value = 'b'
$ws.Range("J410").Value = 'Acknowledge (Backchannel)'
$ws.Range("I417").Value = 'sv'
$ws.Range("J417").Value = 'Statement-opinion'
$ws.Range("I419").Value = 'b'
$ws.Range("J419").Value = 'Acknowledge (Backchannel)'
$ws.Range("I427").Value = 'sv'
$ws.Range("J427").Value = 'Statement-opinion'
$ws.Range("I432").Value = 'aa'
$ws.Range("J432").Value = 'Agree/Accept'
$ws.Range("I440").Value = 'aa'
$ws.Range("J440").Value = 'Agree/Accept'
$ws.Range("I441").Value = 'b'
$ws.Range("J441").Value = 'Acknowledge (Backchannel)'
$ws.Range("I443").Value = 'b'
$ws.Range("J443").Value = 'Acknowledge (Backchannel)'
$ws.Range("I445").Value = 'sd'
$ws.Range("J445").Value = 'Statement-non-opinion'
$ws.Range("I459").Value = 'ba'
$ws.Range("J459").Value = 'Appreciation'
$ws.Range("I529").Value = 'sd'
$ws.Range("J529").Value = 'Statement-non-opinion'
$ws.Range("I536").Value = 'sv'
$ws.Range("J536").Value = 'Statement-opinion'
$ws.Range("I543").Value = 'sd'
$ws.Range("J543").Value = 'Statement-non-opinion'
$ws.Range("I547").Value = 'sd'
$ws.Range("J547").Value = 'Statement-non-opinion'
$ws.Range("I558").Value = 'sd'
$ws.Range("J558").Value = 'Statement-non-opinion'
$ws.Range("I578").Value = 'sv'
$ws.Range("J578").Value = 'Statement-opinion'
$ws.Range("I596").Value = 'aa'
$ws.Range("J596").Value = 'Agree/Accept'
$ws.Range("I598").Value = 'sv'
$ws.Range("J598").Value = 'Statement-opinion'
$ws.Range("I647").Value = 'b'
$ws.Range("J647").Value = 'Acknowledge (Backchannel)'
$ws.Range("I651").Value = 'b'
$ws.Range("J651").Value = 'Acknowledge (Backchannel)'
$ws.Range("I655").Value = 'sd'
$ws.Range("J655").Value = 'Statement-non-opinion'
$ws.Range("I673").Value = 'sd'
$ws.Range("J673").Value = 'Statement-non-opinion'
$ws.Range("I785").Value = 'sd'
$ws.Range("J785").Value = 'Statement-non-opinion'
$ws.Range("I788").Value = 'sd'
$ws.Range("J788").Value = 'Statement-non-opinion'
